$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.57%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.14%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.047"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.57%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08143"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.81%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.082"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8.17%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.866"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.14%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9295"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.12%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1446"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "16.11%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1925"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.45%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09156"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.48%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03443"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.51%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09953"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.24%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006191"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.53%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.837"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.14%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.157"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.31%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.488"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "13.70%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3460"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.65%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1318"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.86%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.825"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.62%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2336"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04392"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.33%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.16%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004202"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-11.19%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.08%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02051"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.32%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05170"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.19%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007470"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.99%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01008"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.13%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1375"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.14%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002131"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.52%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009734"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.07%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006292"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.07%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.81"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.38%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001597"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.71%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
